$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the ticker symbol in A5 (NFTFINSERV25/50 -> NFTFINSERV25_50)
$ws.Range("A5").Value = "NFTFINSERV25_50"

# Update the last saved selection (matches final cursor position in the sheet)
$ws.Range("E23").Select()
